$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4: the customer-id/address data got extra random suffixes appended ---
$ws.Range("D4").Value = "sharonH-830-539-29-142"
$ws.Range("E4").Value = "54 th ave.-930"
$ws.Range("F4").Value = "Lakecity-795"

# --- Row 5: new "create sales team" test case row ---
$ws.Range("A5").Value = "createSalesTeam_ID"
$ws.Range("C5").Value = "moejoe91"
$ws.Range("D5").Value = "Email2-778-764-955"
$ws.Range("E5").Value = "Kevin-541-997-548"
$ws.Range("I5").Value = "clovis"

# --- Row 6: new "create opportunity" test case row ---
$ws.Range("A6").Value = "createOpportunity_ID"
$ws.Range("C6").Value = "moejoe91"
$ws.Range("D6").Value = "100 macs"
$ws.Range("E6").Value = "Kevin"
$ws.Range("F6").Value = 1000

# --- Hyperlinked e-mail cells (added in this order: B5, B6, F5) ---
$ws.Range("B5").Value = "mohitjoe91@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:mohitjoe91@gmail.com")
$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial(-4122)

$ws.Range("B6").Value = "mohitjoe91@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:mohitjoe91@gmail.com")
$ws.Range("B2").Copy()
$ws.Range("B6").PasteSpecial(-4122)

$ws.Range("F5").Value = "Kevin@gmail.com-516-662-668"
$ws.Hyperlinks.Add($ws.Range("F5"), "mailto:Kevin@gmail.com")
$ws.Range("B2").Copy()
$ws.Range("F5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Selection moved to A6 ---
$null = $ws.Range("A6").Select()

Write-Output "sales team created and deleted"
